# Update forecast values in column B (MSTL) for rows 2-49 with the
# newly generated scenario/return statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 125.9144134521484
$ws.Cells.Item(3, 2).Value = 125.9067153930664
$ws.Cells.Item(4, 2).Value = 127.1896591186523
$ws.Cells.Item(5, 2).Value = 127.1970520019531
$ws.Cells.Item(6, 2).Value = 126.2445678710938
$ws.Cells.Item(7, 2).Value = 126.2588882446289
$ws.Cells.Item(8, 2).Value = 126.0774307250977
$ws.Cells.Item(9, 2).Value = 126.0475692749023
$ws.Cells.Item(10, 2).Value = 126.6835403442383
$ws.Cells.Item(11, 2).Value = 126.6560211181641
$ws.Cells.Item(12, 2).Value = 131.4456939697266
$ws.Cells.Item(13, 2).Value = 131.4291534423828
$ws.Cells.Item(14, 2).Value = 145.6224517822266
$ws.Cells.Item(15, 2).Value = 145.6199951171875
$ws.Cells.Item(16, 2).Value = 174.6736450195312
$ws.Cells.Item(17, 2).Value = 174.6851348876953
$ws.Cells.Item(18, 2).Value = 191.0147857666016
$ws.Cells.Item(19, 2).Value = 191.0659637451172
$ws.Cells.Item(20, 2).Value = 202.4267578125
$ws.Cells.Item(21, 2).Value = 202.4864196777344
$ws.Cells.Item(22, 2).Value = 201.2635345458984
$ws.Cells.Item(23, 2).Value = 201.3295135498047
$ws.Cells.Item(24, 2).Value = 192.1592864990234
$ws.Cells.Item(25, 2).Value = 192.2299652099609
$ws.Cells.Item(26, 2).Value = 191.3470764160156
$ws.Cells.Item(27, 2).Value = 191.4209899902344
$ws.Cells.Item(28, 2).Value = 188.7110595703125
$ws.Cells.Item(29, 2).Value = 188.7864532470703
$ws.Cells.Item(30, 2).Value = 191.3690795898438
$ws.Cells.Item(31, 2).Value = 191.444580078125
$ws.Cells.Item(32, 2).Value = 204.1963195800781
$ws.Cells.Item(33, 2).Value = 204.2711334228516
$ws.Cells.Item(34, 2).Value = 230.5118560791016
$ws.Cells.Item(35, 2).Value = 230.5858154296875
$ws.Cells.Item(36, 2).Value = 252.0239715576172
$ws.Cells.Item(37, 2).Value = 252.0974884033203
$ws.Cells.Item(38, 2).Value = 216.789794921875
$ws.Cells.Item(39, 2).Value = 216.8633880615234
$ws.Cells.Item(40, 2).Value = 189.0462188720703
$ws.Cells.Item(41, 2).Value = 189.1197662353516
$ws.Cells.Item(42, 2).Value = 171.1269989013672
$ws.Cells.Item(43, 2).Value = 171.199462890625
$ws.Cells.Item(44, 2).Value = 159.6818695068359
$ws.Cells.Item(45, 2).Value = 159.7518463134766
$ws.Cells.Item(46, 2).Value = 149.7577667236328
$ws.Cells.Item(47, 2).Value = 149.8244934082031
$ws.Cells.Item(48, 2).Value = 150.2505340576172
$ws.Cells.Item(49, 2).Value = 150.3140106201172
